# Applies the commit: "Trade #2 closed at 2026-02-18 10:22:21 - unknown UNKNOWN +0.000%"
#
# Changes:
#  1) Summary sheet        - recalculated aggregate metrics after the close
#  2) Strategy Status sheet - MarketMaking strategy capital / P&L% updated
#  3) All Trades sheet     - trade #4 (row 5) closed out; new OPEN trade #16 appended (row 17)
#  4) MarketMaking sheet   - new OPEN trade #16 appended (row 5)

$wb = $excel.ActiveWorkbook

# Helper: write a date-like text value ("YYYY-MM-DD") into a cell while
# keeping it stored as plain text (matching the rest of the workbook) and
# without leaving a lingering custom number-format style behind.
function Set-TextDate($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# -------------------------------------------------------------------------
# 1) Summary
# -------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.87   # Current Capital
$wsSummary.Range("B4").Value = 0.12      # Total P&L $
$wsSummary.Range("B5").Value = 0.6       # Total P&L %
$wsSummary.Range("B6").Value = 4         # Total Trades
$wsSummary.Range("B8").Value = 1         # Losing Trades
$wsSummary.Range("B9").Value = 50        # Win Rate %

# -------------------------------------------------------------------------
# 2) Strategy Status
# -------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.87      # MarketMaking Capital
$wsStatus.Range("F6").Value = -0.13      # MarketMaking P&L %

# -------------------------------------------------------------------------
# 3) All Trades - close out trade #4 (row 5)
# -------------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("G5").Value = 0.213415
$wsAllTrades.Range("H5").Value = "CLOSED"
$wsAllTrades.Range("I5").Value = -14.6341
$wsAllTrades.Range("J5").Value = -0.18
$wsAllTrades.Range("K5").Value = 99.87
$wsAllTrades.Range("L5").Value = "early_exit"
$wsAllTrades.Range("M5").Value = 0.13

# All Trades - append new OPEN trade #16 as row 17
$wsAllTrades.Range("A17").Value = 16
Set-TextDate $wsAllTrades.Range("B17") "2026-02-18"
$wsAllTrades.Range("C17").Value = "10:22:14"
$wsAllTrades.Range("D17").Value = "MarketMaking"
$wsAllTrades.Range("E17").Value = "DOWN"
$wsAllTrades.Range("F17").Value = 0.25
$wsAllTrades.Range("H17").Value = "OPEN"
$wsAllTrades.Range("I17").Value = 0
$wsAllTrades.Range("J17").Value = 0
$wsAllTrades.Range("K17").Value = 100.05
$wsAllTrades.Range("M17").Value = 0
$wsAllTrades.Range("N17").Value = 0
$wsAllTrades.Range("O17").Value = 0
$wsAllTrades.Range("P17").Value = 0.6
$wsAllTrades.Range("Q17").Value = "Normal spread capture: 202 bps"

# -------------------------------------------------------------------------
# 4) MarketMaking - append the same new OPEN trade #16 as row 5
# -------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A5").Value = 16
Set-TextDate $wsMM.Range("B5") "2026-02-18"
$wsMM.Range("C5").Value = "10:22:14"
$wsMM.Range("D5").Value = "MarketMaking"
$wsMM.Range("E5").Value = "DOWN"
$wsMM.Range("F5").Value = 0.25
$wsMM.Range("H5").Value = "OPEN"
$wsMM.Range("I5").Value = 0
$wsMM.Range("J5").Value = 0
$wsMM.Range("K5").Value = 100.05
$wsMM.Range("L5").Value = 0
$wsMM.Range("M5").Value = 0
$wsMM.Range("N5").Value = 0.6
$wsMM.Range("O5").Value = "Normal spread capture: 202 bps"
$wsMM.Range("Q5").Value = 0
